$wb = $excel.ActiveWorkbook

# --- Sheet: 식당판매 (restaurant sales) ---
$ws1 = $wb.Worksheets.Item("식당판매")
$ws1.Range("C2").Value = 0
$ws1.Range("C3").Value = 0
$ws1.Range("C4").Value = 0
$ws1.Range("C5").Value = 0
$ws1.Range("C6").Value = 0
$ws1.Range("C7").Value = 0
$ws1.Range("C8").Value = 0
$ws1.Range("C9").Value = 0
$ws1.Range("C10").Value = 0
$ws1.Range("C11").Value = 0
$ws1.Range("C12").Value = 0
$ws1.Range("C13").Value = 0

# --- Sheet: 매점판매 (shop sales) ---
$ws2 = $wb.Worksheets.Item("매점판매")
$ws2.Range("C4").Value = 4
$ws2.Range("C5").Value = 21
$ws2.Range("C6").Value = 21

# --- Sheet: 상복 (mourning clothes) ---
$ws4 = $wb.Worksheets.Item("상복")
$ws4.Range("C7").Value = 4

# --- Sheet: 기타 (miscellaneous) ---
$ws5 = $wb.Worksheets.Item("기타")
$ws5.Range("C8").Value = 65
$ws5.Range("C9").Value = 39
$ws5.Range("C10").Value = 223
$ws5.Range("C11").Value = 44

# --- Sheet: 세트 (set) ---
$ws6 = $wb.Worksheets.Item("세트")
$ws6.Range("C1").Value = "사용 수량"
$ws6.Range("E1").Value = ""

$ws6.Range("A2").Value = "수저"
$ws6.Range("B2").Value = 780
$ws6.Range("C2").Value = 3

$ws6.Range("A3").Value = "대패삼겹살"
$ws6.Range("B3").Value = 2000
$ws6.Range("C3").Value = 1

$ws6.Range("A4").Value = "치즈김밥"
$ws6.Range("B4").Value = 2500
$ws6.Range("C4").Value = 2
